$wb = $excel.ActiveWorkbook

# Update the "Status" text from "Ready for handoff" to "In Translation" on every sheet
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ("$val" -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# The shorter replacement text narrows the "Status" columns (report regenerated/re-fit).
# Set the resulting column widths to match the regenerated report.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.45
$wsOverview.Columns.Item(6).ColumnWidth = 12.45

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.45

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.45
